$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value changes per diff (Jogos_da_Semana_FlashScore_2025-05-13.xlsx)
$ws.Range("G5").Value = 2.35
$ws.Range("I5").Value = 3
$ws.Range("U5").Value = 11
$ws.Range("AE5").Value = 9.5
$ws.Range("N6").Value = 2.04
$ws.Range("O6").Value = 1.86
$ws.Range("G12").Value = 2.72
$ws.Range("H12").Value = 2.75
$ws.Range("I12").Value = 2.8
$ws.Range("M12").Value = 2.27
$ws.Range("N12").Value = 2.4
$ws.Range("Q12").Value = 2.2
$ws.Range("R12").Value = 1.98
$ws.Range("T12").Value = 6.7
$ws.Range("U12").Value = 12.5
$ws.Range("AA12").Value = 5.5
$ws.Range("AB12").Value = 16.5
$ws.Range("AC12").Value = 100
$ws.Range("AG12").Value = 10.75
$ws.Range("AJ12").Value = 45
$ws.Range("G13").Value = 2.32
$ws.Range("H13").Value = 2.95
$ws.Range("I13").Value = 3.1
$ws.Range("L13").Value = 1.45
$ws.Range("M13").Value = 2.37
$ws.Range("N13").Value = 2.32
$ws.Range("O13").Value = 1.47
$ws.Range("P13").Value = 1.5
$ws.Range("Q13").Value = 2.27
$ws.Range("R13").Value = 1.98
$ws.Range("S13").Value = 1.65
$ws.Range("T13").Value = 6.1
$ws.Range("U13").Value = 10
$ws.Range("V13").Value = 9.75
$ws.Range("W13").Value = 24
$ws.Range("X13").Value = 23
$ws.Range("Y13").Value = 40
$ws.Range("Z13").Value = 6.8
$ws.Range("AA13").Value = 5.9
$ws.Range("AB13").Value = 17.5
$ws.Range("AC13").Value = 110
$ws.Range("AE13").Value = 7.5
$ws.Range("AF13").Value = 14.5
$ws.Range("AG13").Value = 11.5
$ws.Range("AH13").Value = 40
$ws.Range("AI13").Value = 32
$ws.Range("AJ13").Value = 50
$ws.Range("N23").Value = 2.25
$ws.Range("O23").Value = 1.62
$ws.Range("L24").Value = 1.18
$ws.Range("M24").Value = 4.5
$ws.Range("G25").Value = 3.75
$ws.Range("H25").Value = 3.1
$ws.Range("I25").Value = 2.1
$ws.Range("J25").Value = 1.06
$ws.Range("K25").Value = 10
$ws.Range("U25").Value = 19
$ws.Range("Y25").Value = 34
$ws.Range("AF25").Value = 10
$ws.Range("AH25").Value = 19
$ws.Range("AI25").Value = 17
$ws.Range("G26").Value = 1.62
$ws.Range("W26").Value = 12
$ws.Range("AD26").Value = 251
$ws.Range("AF26").Value = 29
$ws.Range("G27").Value = 3.25
$ws.Range("H27").Value = 3.6
$ws.Range("I27").Value = 2.1
$ws.Range("K27").Value = 12
$ws.Range("U27").Value = 17
$ws.Range("W27").Value = 34
$ws.Range("Y27").Value = 29
$ws.Range("AB27").Value = 13
$ws.Range("AE27").Value = 8.5
$ws.Range("AF27").Value = 11
$ws.Range("AG27").Value = 9
$ws.Range("AH27").Value = 19
$ws.Range("AI27").Value = 17
$ws.Range("N29").Value = 1.6
$ws.Range("O29").Value = 2.3
$ws.Range("N30").Value = 2.2
$ws.Range("O30").Value = 1.65
$ws.Range("I31").Value = 4.3
$ws.Range("AC31").Value = 90
$ws.Range("AG31").Value = 14.5
$ws.Range("I33").Value = 3.8
$ws.Range("P33").Value = 1.55
$ws.Range("Q33").Value = 2.35
$ws.Range("G34").Value = 3.8
$ws.Range("H34").Value = 3.3
$ws.Range("I34").Value = 2
$ws.Range("AB34").Value = 15
$ws.Range("AD34").Value = 301
$ws.Range("AE34").Value = 7
$ws.Range("AF34").Value = 9
$ws.Range("G42").Value = 3.3
$ws.Range("H42").Value = 3.75
$ws.Range("I42").Value = 1.95
$ws.Range("K42").Value = 8.75
$ws.Range("M42").Value = 4.05
$ws.Range("N42").Value = 1.6
$ws.Range("O42").Value = 2.2
$ws.Range("T42").Value = 13
$ws.Range("U42").Value = 20
$ws.Range("Y42").Value = 28
$ws.Range("Z42").Value = 8.75
$ws.Range("AA42").Value = 7.5
$ws.Range("AB42").Value = 12.5
$ws.Range("AE42").Value = 9.75
$ws.Range("AF42").Value = 11
$ws.Range("AI42").Value = 14
$ws.Range("AJ42").Value = 20
